# Switch the presentation's theme colour scheme from "Литейная" (Foundry)
# to "Красный и оранжевый" (Red Orange).
#
# The OOXML theme (ppt/theme/theme1.xml) keeps its 12 colour slots in a
# fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. The
# PowerPoint object model exposes that same ordering as a 1-based index
# through ThemeColorScheme.Colors(i).RGB, where .RGB uses the classic
# VBA encoding (R + G*256 + B*65536). Converting each target hex colour
# to that encoding and writing it back reproduces the diff exactly.

function HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# All slides share the single slide master's theme.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Best-effort: rename the colour scheme itself to match the new palette.
try { $tcs.Name = "Красный и оранжевый" } catch { }

# Index map (1-based) -> new sRGB hex value.
# 1 dk1       000000  (unchanged)
# 2 lt1       FFFFFF  (unchanged)
# 3 dk2       505046
# 4 lt2       EEECE1
# 5 accent1   E84C22
# 6 accent2   FFBD47
# 7 accent3   B64926
# 8 accent4   FF8427
# 9 accent5   CC9900
# 10 accent6  B22600
# 11 hlink    CC9900
# 12 folHlink 666699

$tcs.Colors(3).RGB = HexToVbaRgb "505046"
$tcs.Colors(4).RGB = HexToVbaRgb "EEECE1"
$tcs.Colors(5).RGB = HexToVbaRgb "E84C22"
$tcs.Colors(6).RGB = HexToVbaRgb "FFBD47"
$tcs.Colors(7).RGB = HexToVbaRgb "B64926"
$tcs.Colors(8).RGB = HexToVbaRgb "FF8427"
$tcs.Colors(9).RGB = HexToVbaRgb "CC9900"
$tcs.Colors(10).RGB = HexToVbaRgb "B22600"
$tcs.Colors(11).RGB = HexToVbaRgb "CC9900"
$tcs.Colors(12).RGB = HexToVbaRgb "666699"
